$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the data set. It belongs right
# after the existing row 63 (row 64 in the original sheet), so insert a
# fresh row there; Excel shifts every row from 64..142 down to 65..143
# and carries the column D date-number formatting (style) along with it.
$ws.Rows("64:64").Insert()

# Populate the newly inserted row 64 with the new record's data.
$ws.Cells.Item(64, 1).Value = 1
$ws.Cells.Item(64, 2).Value = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(64, 3).Value = 'Arica y Parinacota'
$ws.Cells.Item(64, 4).Value = 45100
$ws.Cells.Item(64, 5).Value = 15
$ws.Cells.Item(64, 6).Value = 100112021
$ws.Cells.Item(64, 7).Value = 'Ají'
$ws.Cells.Item(64, 8).Value = 'Inferno'
$ws.Cells.Item(64, 9).Value = 'Primera'
$ws.Cells.Item(64, 10).Value = 150
$ws.Cells.Item(64, 11).Value = 10000
$ws.Cells.Item(64, 12).Value = 11000
$ws.Cells.Item(64, 13).Value = 10500
$ws.Cells.Item(64, 14).Value = '$/caja 15 kilos'
$ws.Cells.Item(64, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(64, 16).Value = 700
$ws.Cells.Item(64, 17).Value = 15
$ws.Cells.Item(64, 18).Value = 'Hortaliza'
